$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "John T. Hawkins"
$ws.Range("B2").Value = "JohnTHawkins@jourrapide.com            This is a real email address. Click here to activate it!"
$ws.Range("C2").Value = "`n                                        3538 Filbert StreetChester, PA 19013                                        "
$ws.Range("D2").Value = "610-876-7584"

$ws.Range("A3").Value = "Beverly J. Chambers"
$ws.Range("B3").Value = "BeverlyJChambers@teleworm.us            This is a real email address. Click here to activate it!"
$ws.Range("C3").Value = "`n                                        4055 Pallet StreetWest Nyack, NY 10994                                        "
$ws.Range("D3").Value = "914-346-4627"

$ws.Range("A4").Value = "Gordon P. Fields"
$ws.Range("B4").Value = "GordonPFields@teleworm.us            This is a real email address. Click here to activate it!"
$ws.Range("C4").Value = "`n                                        4402 Modoc AlleyDixie, ID 83525                                        "
$ws.Range("D4").Value = "208-842-0461"

$ws.Rows("3").RowHeight = 99.75
$ws.Rows("4").RowHeight = 85.5
